# Atualização de bases das ligas, do dia: 27-04-2024 às 09:20
#
# A new finished match is inserted at row 176 (id 174) and the matches that
# used to occupy rows 176-180 shift down to rows 177-181 (their "live" odds
# also get refreshed at the same time). Row 181 is a brand new row holding
# what used to be row 180's match (id 6885578), with refreshed odds.
#
# Columns B:F (match id / league / date / home / away) are moved down one
# row at a time using Range.Copy, bottom-up, so the cell types (several of
# the "id" values in column B are stored as shared-string text, not
# numbers) and styles (bold border on column A, date format on column D)
# are preserved exactly instead of being coerced by a plain .Value
# assignment. Column A (the simple running "id" index 174..179) is NOT
# copied - it is left alone / set explicitly afterwards, since it must stay
# in row order rather than shifting down with the match data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B180:F180").Copy($ws.Range("B181:F181"))
$ws.Range("B179:F179").Copy($ws.Range("B180:F180"))
$ws.Range("B178:F178").Copy($ws.Range("B179:F179"))
$ws.Range("B177:F177").Copy($ws.Range("B178:F178"))
$ws.Range("B176:F176").Copy($ws.Range("B177:F177"))

# --- Row 181 (brand new row): id 179, refreshed odds for what used to be
#     row 180's match.
$ws.Range("A181").Value = 179
$ws.Range("J181").Value = 4.333
$ws.Range("K181").Value = 3.75
$ws.Range("L181").Value = 1.666
$ws.Range("M181").Value = 5.75
$ws.Range("N181").Value = 4.2
$ws.Range("O181").Value = 1.5
$ws.Range("P181").Value = 1
$ws.Range("Q181").Value = 2.025
$ws.Range("R181").Value = 1.825
$ws.Range("S181").Value = 3
$ws.Range("T181").Value = 1.975
$ws.Range("U181").Value = 1.875
$ws.Range("V181").Value = 0
$ws.Range("W181").Value = 0
$ws.Range("X181").Value = 0

# --- Row 180: now what used to be row 179's match, refreshed odds.
$ws.Range("J180").Value = 2.6
$ws.Range("K180").Value = 3.25
$ws.Range("L180").Value = 2.5
$ws.Range("M180").Value = 2.9
$ws.Range("N180").Value = 3.25
$ws.Range("O180").Value = 2.3
$ws.Range("P180").Value = 0.25
$ws.Range("Q180").Value = 1.825
$ws.Range("R180").Value = 2.025
$ws.Range("S180").Value = 2.5
$ws.Range("T180").Value = 1.925
$ws.Range("U180").Value = 1.925

# --- Row 179: now what used to be row 178's match, refreshed odds.
$ws.Range("J179").Value = 1.75
$ws.Range("K179").Value = 3.6
$ws.Range("L179").Value = 4
$ws.Range("M179").Value = 1.75
$ws.Range("N179").Value = 3.75
$ws.Range("O179").Value = 4.2
$ws.Range("P179").Value = -0.75
$ws.Range("Q179").Value = 2.025
$ws.Range("R179").Value = 1.825
$ws.Range("S179").Value = 2.75
$ws.Range("T179").Value = 1.95
$ws.Range("U179").Value = 1.9

# --- Row 178: now what used to be row 177's match, refreshed odds.
$ws.Range("J178").Value = 2
$ws.Range("K178").Value = 3.3
$ws.Range("L178").Value = 3.4
$ws.Range("M178").Value = 1.85
$ws.Range("N178").Value = 3.2
$ws.Range("O178").Value = 4.333
$ws.Range("P178").Value = -0.5
$ws.Range("Q178").Value = 1.925
$ws.Range("R178").Value = 1.925
$ws.Range("S178").Value = 2.5
$ws.Range("T178").Value = 2
$ws.Range("U178").Value = 1.85

# --- Row 177: now what used to be row 176's match, refreshed odds.
$ws.Range("J177").Value = 1.4
$ws.Range("K177").Value = 4.5
$ws.Range("L177").Value = 6.5
$ws.Range("M177").Value = 1.25
$ws.Range("N177").Value = 5.5
$ws.Range("O177").Value = 10
$ws.Range("P177").Value = -1.75
$ws.Range("Q177").Value = 1.975
$ws.Range("R177").Value = 1.875
$ws.Range("S177").Value = 3
$ws.Range("T177").Value = 1.8
$ws.Range("U177").Value = 2.05

# --- Row 176: becomes a brand-new, already-finished match.
$ws.Range("B176").Value = 6943666
$ws.Range("C176").Value = "Hungary NB I"
$ws.Range("D176").Value = 45408.625
$ws.Range("E176").Value = "Zalaegerszegi TE"
$ws.Range("F176").Value = "Puskas Academy"
$ws.Range("G176").Value = 1
$ws.Range("H176").Value = 0
$ws.Range("I176").Value = "H"
$ws.Range("J176").Value = 3.6
$ws.Range("K176").Value = 3.4
$ws.Range("L176").Value = 1.909
$ws.Range("M176").Value = 3.75
$ws.Range("N176").Value = 3.6
$ws.Range("O176").Value = 1.85
$ws.Range("P176").Value = 0.5
$ws.Range("Q176").Value = 1.95
$ws.Range("R176").Value = 1.9
$ws.Range("S176").Value = 2.75
$ws.Range("T176").Value = 1.9
$ws.Range("U176").Value = 1.95
$ws.Range("V176").Value = 2.75
$ws.Range("W176").Value = -1
$ws.Range("X176").Value = -1
$ws.Range("Y176").Value = 0.95
$ws.Range("Z176").Value = -1
$ws.Range("AA176").Value = -1
$ws.Range("AB176").Value = 0.95

Write-Host "Rows 176-181 updated"
